$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells I1:W1 (Round_6 .. Round_20), matching style/format of existing headers ---
$ws.Range("H1").Copy()
$ws.Range("I1:W1").PasteSpecial(-4122)

$roundNum = 6
for ($col = 9; $col -le 23; $col++) {
    $ws.Cells.Item(1, $col).Value = "Round_" + $roundNum
    $roundNum++
}

# --- Row 2: rename participant/team, update round scores, add new rounds ---
$ws.Range("A2").Value = "Yamai"
$ws.Range("C2").Value = "Yamai Syndicate"

$row2Values = @(6,4,6,2,4,6,7,6,10,7,3,1,0,6,3,8,8,6,1,0)
for ($i = 0; $i -lt $row2Values.Length; $i++) {
    $ws.Cells.Item(2, 4 + $i).Value = $row2Values[$i]
}

# --- Row 3: rename participant/team, update round scores, add new rounds ---
$ws.Range("A3").Value = "Kalani"
$ws.Range("C3").Value = "Yamai Syndicate"

$row3Values = @(6,9,0,6,7,7,9,4,1,6,7,10,9,3,10,3,0,6,7,4)
for ($i = 0; $i -lt $row3Values.Length; $i++) {
    $ws.Cells.Item(3, 4 + $i).Value = $row3Values[$i]
}

# --- Remove rows 4, 5, 6 (Cyber, Nozomi, Inizio) ---
$ws.Range("A4:W6").EntireRow.Delete()
